$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# wdHeaderFooterPrimary = 1, wdHeaderFooterFirstPage = 2
$hdrPrimary = $sec.Headers.Item(1)
$hdrFirst   = $sec.Headers.Item(2)
$ftrPrimary = $sec.Footers.Item(1)
$ftrFirst   = $sec.Footers.Item(2)

# Headers contain the BTec_Logo-Orange picture: image1.jpg -> image2.jpg
if ($hdrPrimary.Exists -and $hdrPrimary.Range.InlineShapes.Count -gt 0) {
    $hdrPrimary.Range.InlineShapes.Item(1).Name = "image2.jpg"
}
if ($hdrFirst.Exists -and $hdrFirst.Range.InlineShapes.Count -gt 0) {
    $hdrFirst.Range.InlineShapes.Item(1).Name = "image2.jpg"
}

# Footers contain the PearsonLogo picture: image2.png -> image1.png
if ($ftrPrimary.Exists -and $ftrPrimary.Range.InlineShapes.Count -gt 0) {
    $ftrPrimary.Range.InlineShapes.Item(1).Name = "image1.png"
}
if ($ftrFirst.Exists -and $ftrFirst.Range.InlineShapes.Count -gt 0) {
    $ftrFirst.Range.InlineShapes.Item(1).Name = "image1.png"
}
